# Combine section#goals task row removed (responsive table made for skills section).
# Delete the task row containing "Combine section#goals unordered lists into a
# single responsive table" (row 3), which shifts the subsequent rows up by one
# and automatically removes the now-unused shared string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A3:C3").EntireRow.Delete() | Out-Null

# Update the selection to match the post-edit cursor position.
$ws.Range("B8").Select() | Out-Null
